$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 179, shifting rows 179:265 down to 180:266
$ws.Rows.Item(179).Insert()

# Fill the new row 179 with its values
$ws.Cells.Item(179, 1).Value = 3
$ws.Cells.Item(179, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(179, 3).Value = "Coquimbo"
$ws.Cells.Item(179, 4).Value = 44992
$ws.Cells.Item(179, 5).Value = 5
$ws.Cells.Item(179, 6).Value = 100112030
$ws.Cells.Item(179, 7).Value = "Poroto granado"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 45
$ws.Cells.Item(179, 11).Value = 30000
$ws.Cells.Item(179, 12).Value = 30000
$ws.Cells.Item(179, 13).Value = 30000
$ws.Cells.Item(179, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(179, 15).Value = "Provincia de Petorca"
$ws.Cells.Item(179, 16).Value = 1200
$ws.Cells.Item(179, 17).Value = 25
$ws.Cells.Item(179, 18).Value = "Hortaliza"
